# Update Código Efecto Hall
# The VL columns (F:I) were measured in mV; the units are corrected to V
# and the VL_err column (J) values are rescaled from 0.1 (mV) to 1E-3 (V).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 1 headers: change unit label from "mV" to "V"
$ws.Range("F1").Value = "VL_-20mA (V)"
$ws.Range("G1").Value = "VL_-10mA (V)"
$ws.Range("H1").Value = "VL_10mA (V)"
$ws.Range("I1").Value = "VL_20mA (V)"
$ws.Range("J1").Value = "VL_err (V)"

# Column J (VL_err), rows 2-16: 0.1 mV -> 1E-3 V
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 10).Value = 0.001
}

# Reflect the recorded selection change on the sheet (G11 was last selected).
$ws.Range("G11").Select()
